$wb = $excel.ActiveWorkbook

# Scheduled runner update: refresh market-board derived leve profit figures
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 75000
$ws.Range("J81").Value = 75000
$ws.Range("L81").Value = 75000
$ws.Range("N81").Value = -76996

$ws.Range("H84").Value = 75000
$ws.Range("J84").Value = 75000
$ws.Range("L84").Value = 225000
$ws.Range("N84").Value = -234984

$ws.Range("H98").Value = 2577.4614
$ws.Range("I98").Value = 2650.96
$ws.Range("J98").Value = 740
$ws.Range("K98").Value = 2650.96
$ws.Range("L98").Value = 740
$ws.Range("M98").Value = -1152.96
$ws.Range("N98").Value = -3736

$ws.Range("H112").Value = 2003.92
$ws.Range("I112").Value = 1066.6666
$ws.Range("J112").Value = 2131.7273
$ws.Range("K112").Value = 3199.9998
$ws.Range("L112").Value = 6395.1819
$ws.Range("M112").Value = -2091.9998
$ws.Range("N112").Value = -8611.1819

$ws.Range("H122").Value = 2577.4614
$ws.Range("I122").Value = 2650.96
$ws.Range("J122").Value = 740
$ws.Range("K122").Value = 7952.88
$ws.Range("L122").Value = 2220
$ws.Range("M122").Value = -5502.88
$ws.Range("N122").Value = -7120

$ws.Range("H137").Value = 1477558.5
$ws.Range("J137").Value = 8003.069
$ws.Range("L137").Value = 24009.207
$ws.Range("N137").Value = -29109.207

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 27789.688
$ws.Range("I45").Value = 31259.5
$ws.Range("J45").Value = 3501
$ws.Range("K45").Value = 31259.5
$ws.Range("L45").Value = 3501
$ws.Range("M45").Value = -30882.5
$ws.Range("N45").Value = -4255

$ws.Range("H61").Value = 6084.1665
$ws.Range("I61").Value = 2990
$ws.Range("K61").Value = 2990
$ws.Range("M61").Value = -2778

$ws.Range("H136").Value = 6084.1665
$ws.Range("I136").Value = 2990
$ws.Range("K136").Value = 8970
$ws.Range("M136").Value = -6420

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 7200.933
$ws.Range("I99").Value = 3692.875
$ws.Range("J99").Value = 11210.143
$ws.Range("K99").Value = 3692.875
$ws.Range("L99").Value = 11210.143
$ws.Range("M99").Value = -2194.875
$ws.Range("N99").Value = -14206.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1414.0714
$ws.Range("I16").Value = 1252.4762
$ws.Range("J16").Value = 1898.8572
$ws.Range("K16").Value = 1252.4762
$ws.Range("L16").Value = 1898.8572
$ws.Range("M16").Value = -965.4762000000001
$ws.Range("N16").Value = -2472.8572

$ws.Range("H31").Value = 5701.28
$ws.Range("I31").Value = 4064.1428
$ws.Range("J31").Value = 7784.909
$ws.Range("K31").Value = 4064.1428
$ws.Range("L31").Value = 7784.909
$ws.Range("M31").Value = -3769.1428
$ws.Range("N31").Value = -8374.909

$ws.Range("H34").Value = 5701.28
$ws.Range("I34").Value = 4064.1428
$ws.Range("J34").Value = 7784.909
$ws.Range("K34").Value = 4064.1428
$ws.Range("L34").Value = 7784.909
$ws.Range("M34").Value = -3862.1428
$ws.Range("N34").Value = -8188.909

$ws.Range("H51").Value = 58157.25
$ws.Range("J51").Value = 57528.668
$ws.Range("L51").Value = 57528.668
$ws.Range("N51").Value = -59000.668

$ws.Range("H58").Value = 2274.762
$ws.Range("I58").Value = 1181.5
$ws.Range("J58").Value = 3268.6365
$ws.Range("K58").Value = 1181.5
$ws.Range("L58").Value = 3268.6365
$ws.Range("M58").Value = -978.5
$ws.Range("N58").Value = -3674.6365

$ws.Range("H59").Value = 74813
$ws.Range("J59").Value = 79750.664
$ws.Range("L59").Value = 79750.664
$ws.Range("N59").Value = -82040.664

$ws.Range("H60").Value = 12700.429
$ws.Range("J60").Value = 12700.429
$ws.Range("L60").Value = 12700.429
$ws.Range("N60").Value = -13722.429

$ws.Range("H61").Value = 58157.25
$ws.Range("J61").Value = 57528.668
$ws.Range("L61").Value = 57528.668
$ws.Range("N61").Value = -58224.668

$ws.Range("H68").Value = 86673.25
$ws.Range("J68").Value = 86673.25
$ws.Range("L68").Value = 86673.25
$ws.Range("N68").Value = -88171.25

$ws.Range("H71").Value = 86673.25
$ws.Range("J71").Value = 86673.25
$ws.Range("L71").Value = 260019.75
$ws.Range("N71").Value = -267507.75

$ws.Range("H74").Value = 55381.5
$ws.Range("J74").Value = 55381.5
$ws.Range("L74").Value = 55381.5
$ws.Range("N74").Value = -57129.5

$ws.Range("H77").Value = 55381.5
$ws.Range("J77").Value = 55381.5
$ws.Range("L77").Value = 166144.5
$ws.Range("N77").Value = -174880.5

$ws.Range("H86").Value = 4077
$ws.Range("J86").Value = 3936
$ws.Range("L86").Value = 3936
$ws.Range("N86").Value = -6182

$ws.Range("H88").Value = 37875
$ws.Range("I88").Value = 10000
$ws.Range("K88").Value = 10000
$ws.Range("M88").Value = -9594

$ws.Range("H89").Value = 4077
$ws.Range("J89").Value = 3936
$ws.Range("L89").Value = 19680
$ws.Range("N89").Value = -30912

$ws.Range("H91").Value = 37875
$ws.Range("I91").Value = 10000
$ws.Range("K91").Value = 10000
$ws.Range("M91").Value = -8596

$ws.Range("H92").Value = 40871.5
$ws.Range("J92").Value = 40871.5
$ws.Range("L92").Value = 40871.5
$ws.Range("N92").Value = -45863.5

$ws.Range("H113").Value = 1414.0714
$ws.Range("I113").Value = 1252.4762
$ws.Range("J113").Value = 1898.8572
$ws.Range("K113").Value = 1252.4762
$ws.Range("L113").Value = 1898.8572
$ws.Range("M113").Value = 917.5237999999999
$ws.Range("N113").Value = -6238.8572

$ws.Range("H117").Value = 60000
$ws.Range("J117").Value = 60000
$ws.Range("L117").Value = 60000
$ws.Range("N117").Value = -69178

$ws.Range("H132").Value = 3999.3333
$ws.Range("I132").Value = 3999.3333
$ws.Range("K132").Value = 11997.9999
$ws.Range("M132").Value = -9467.999899999999

$ws.Range("H136").Value = 2274.762
$ws.Range("I136").Value = 1181.5
$ws.Range("J136").Value = 3268.6365
$ws.Range("K136").Value = 3544.5
$ws.Range("L136").Value = 9805.9095
$ws.Range("M136").Value = -994.5
$ws.Range("N136").Value = -14905.9095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 31586644
$ws.Range("J4").Value = 7803749.5
$ws.Range("L4").Value = 23411248.5
$ws.Range("N4").Value = -23411472.5

$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 58825640
$ws.Range("I80").Value = 90910930
$ws.Range("K80").Value = 90910930
$ws.Range("M80").Value = -90909932

$ws.Range("H83").Value = 58825640
$ws.Range("I83").Value = 90910930
$ws.Range("K83").Value = 454554650
$ws.Range("M83").Value = -454549658

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3073.1538
$ws.Range("I68").Value = 3169.25
$ws.Range("J68").Value = 2919.4
$ws.Range("K68").Value = 3169.25
$ws.Range("L68").Value = 2919.4
$ws.Range("M68").Value = -2420.25
$ws.Range("N68").Value = -4417.4

$ws.Range("H71").Value = 3073.1538
$ws.Range("I71").Value = 3169.25
$ws.Range("J71").Value = 2919.4
$ws.Range("K71").Value = 15846.25
$ws.Range("L71").Value = 14597
$ws.Range("M71").Value = -12102.25
$ws.Range("N71").Value = -22085

Write-Output "edit.ps1 applied"
